# Generate Report for Handoff
# Update the handoff/handback timestamps for the
# "9d80176f-0260-46fc-80c5-a33fedd9ee90" entry (row 5) across the
# Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-01-17 11:01:13"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-17 11:01:10"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-17 11:01:13"
